$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.435.77"
$ws.Range("E2").Value = "  -2.41%  "

$ws.Range("D3").Value = "2.293.46"
$ws.Range("E3").Value = "  -3.45%  "

$ws.Range("E4").Value = "  -0.08%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "537.33"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "127.52"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -5.03%  "

$ws.Range("E7").Value = "  -0.11%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.564"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -4.31%  "

$ws.Range("D9").Value = "2.292.54"
$ws.Range("E9").Value = "  -3.35%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0998"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.46"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -2.00%  "

$ws.Range("E12").Value = "  -1.11%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.328"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.81%  "

$ws.Range("D14").Value = "2.705.45"
$ws.Range("E14").Value = "  -3.56%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "59.308.67"
$ws.Range("E15").Value = "  -2.55%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "22.91"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -5.91%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.0000130"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -3.29%  "

$ws.Range("D18").Value = "2.298.65"
$ws.Range("E18").Value = "  -3.67%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "10.31"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -4.79%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "3.99"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -6.02%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "307.40"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -3.70%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.44"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -7.36%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "62.66"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.71%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.167"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -3.78%  "

$ws.Range("E26").Value = "  +0.12%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "7.63"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -6.70%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.32"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.27%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "171.23"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.17"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("E31").Value = "  -3.94%  "

$ws.Range("D32").Value = "0.0₃0705"
$ws.Range("E32").Value = "  -6.42%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.72"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -4.34%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.374"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.59%  "

$ws.Range("E36").Value = "  -7.20%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "17.59"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.04%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.94"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -6.74%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "307.01"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -6.08%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "37.53"
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.48"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -6.16%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "134.66"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -8.15%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.38"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -3.87%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0931"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.83%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.562"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.92%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.0485"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.86%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "18.28"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -6.77%  "

$ws.Range("D49").Value = "0.0₆0216"
$ws.Range("E49").Value = "  +18.77%  "

$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("E51").Value = "  -0.57%  "

